$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "W-" data rows (rows 8-13), mirroring the structure of the existing
# "W+" rows (2-7) but with updated eta/value/stat_u/syst numbers.
# NOTE: scientific-notation literals (e.g. 5.1E-2) are not supported by the
# PowerShell parser here, so all numbers are written in plain decimal form.

$data = @(
    @{ row=8;  D=-2;   E=-1.1; F=-1.27; I=0.241; J=0.146; K=0.01  },
    @{ row=9;  D=-1.1; E=-0.5; F=-0.74; I=0.26;  J=0.051; K=0.01  },
    @{ row=10; D=-0.5; E=0;    F=-0.27; I=0.281; J=0.056; K=0.011 },
    @{ row=11; D=0;    E=0.5;  F=0.27;  I=0.239; J=0.056; K=0.01  },
    @{ row=12; D=0.5;  E=1.1;  F=0.74;  I=0.385; J=0.051; K=0.014 },
    @{ row=13; D=1.1;  E=2;    F=1.27;  I=0.205; J=0.148; K=0.009 }
)

foreach ($d in $data) {
    $r = $d.row

    $ws.Cells.Item($r, 1).Value = "STAR"     # A col
    $ws.Cells.Item($r, 2).Value = "pp"       # B target
    $ws.Cells.Item($r, 3).Value = 510        # C cms
    $ws.Cells.Item($r, 4).Value = $d.D       # D eta_min
    $ws.Cells.Item($r, 5).Value = $d.E       # E eta_max
    $ws.Cells.Item($r, 6).Value = $d.F       # F eta
    $ws.Cells.Item($r, 7).Value = 25         # G pt_min
    $ws.Cells.Item($r, 8).Value = "W-"       # H boson
    $ws.Cells.Item($r, 9).Value = $d.I       # I value
    $ws.Cells.Item($r, 10).Value = $d.J      # J stat_u
    $ws.Cells.Item($r, 11).Value = $d.K      # K syst
    $ws.Cells.Item($r, 12).Value = "A_L"     # L obs
    $ws.Cells.Item($r, 13).Value = "eta"     # M diff
    $ws.Cells.Item($r, 14).Value = 0.005     # N lumi

    $ws.Range("A$r`:N$r").HorizontalAlignment = -4108   # xlCenter, matches style s="1"
}

$ws.Range("O8").Formula = "=0.033*I8"
$ws.Range("O9:O13").Formula = "=0.033*I9"
$ws.Range("O8:O13").HorizontalAlignment = -4108

$ws.Range("P8").Formula = "=SQRT(K8*K8-N8*N8-O8*O8)"
$ws.Range("P9:P13").Formula = "=SQRT(K9*K9-N9*N9-O9*O9)"

$ws.Range("Q8").Formula = "=SQRT(N8*N8+O8*O8)"
$ws.Range("Q9:Q13").Formula = "=SQRT(N9*N9+O9*O9)"

$ws.Range("N15").Select()
